$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.418.84'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.808.27'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.37'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.44'
$ws.Range("E6").Value = '  +1.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.807.03'
$ws.Range("E7").Value = '  +1.68%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  +0.19%  '

$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  -0.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("E13").Value = '  -0.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.26'
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.429.22'
$ws.Range("E15").Value = '  +1.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.786.76'
$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.434.85'
$ws.Range("E17").Value = '  +0.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.03'
$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.00'
$ws.Range("E19").Value = '  -0.80%  '

$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.112'
$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.00'
$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '466.04'
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.702'
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("E24").Value = '  +9.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.08'
$ws.Range("E25").Value = '  +1.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.18'
$ws.Range("E26").Value = '  -2.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.96'
$ws.Range("E27").Value = '  -1.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.15'
$ws.Range("E28").Value = '  -0.23%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("E30").Value = '  -0.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.34'
$ws.Range("E31").Value = '  -0.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.28'
$ws.Range("E32").Value = '  +1.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("E33").Value = '  -2.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.18'
$ws.Range("E34").Value = '  +0.67%  '

$ws.Range("E35").Value = '  +0.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.749.13'
$ws.Range("E36").Value = '  +1.34%  '

$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.45'
$ws.Range("E38").Value = '  +0.61%  '

$ws.Range("E39").Value = '  +1.02%  '

$ws.Range("E40").Value = '  +1.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.80'
$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.29'
$ws.Range("E44").Value = '  +17.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.302'
$ws.Range("E45").Value = '  -1.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.16'
$ws.Range("E46").Value = '  +3.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.92'
$ws.Range("E47").Value = '  -0.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.44'
$ws.Range("E48").Value = '  -1.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '396.90'
$ws.Range("E49").Value = '  +0.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '146.37'
$ws.Range("E50").Value = '  +1.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.806.12'
$ws.Range("E51").Value = '  +4.58%  '
